$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (cells K2:T2)
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.452919
$ws.Range("N2").Value = 0.905838
$ws.Range("O2").Value = 0.03428003430836644
$ws.Range("P2").Value = 0.03049201402277582
$ws.Range("Q2").Value = 0.0549089555865
$ws.Range("R2").Value = 0.219635822346
$ws.Range("S2").Value = 0.03428003430836644
$ws.Range("T2").Value = 0.03049201402277582

# Update row 3 values (cells M3:T3)
$ws.Range("M3").Value = 3.282732333333333
$ws.Range("N3").Value = 9.848196999999999
$ws.Range("O3").Value = 0.2484598283839946
$ws.Range("P3").Value = 0.3315066943791922
$ws.Range("Q3").Value = 0.3979771303331666
$ws.Range("R3").Value = 2.387862781998999
$ws.Range("S3").Value = 0.2484598283839946
$ws.Range("T3").Value = 0.3315066943791922

# Update row 4 values (cells M4:T4)
$ws.Range("M4").Value = 9.476675
$ws.Range("N4").Value = 18.95335
$ws.Range("O4").Value = 0.7172601373076389
$ws.Range("P4").Value = 0.638001291598032
$ws.Range("Q4").Value = 1.1488904786125
$ws.Range("R4").Value = 4.59556191445
$ws.Range("S4").Value = 0.7172601373076389
$ws.Range("T4").Value = 0.638001291598032

# Remove row 5 (Neutrophils target-cluster row) entirely
$ws.Rows.Item(5).Delete()
